# Reverse changes to ENA raw sequencing reads:
#  - Version 1.0.5 -> 1.0.4
#  - Insert "Characteristic [Sequencing Library Name]" / Term Source REF / Term
#    Accession Number (NCIT:C182058) columns before "Characteristic [Sequencing
#    Library Source Indicator]"
#  - Insert "Characteristic [MD5 Checksum]" / Term Source REF / Term Accession
#    Number (NCIT:C171276) columns before "Output [Data]", and drop the old
#    "Data Format" / "Data Selector Format" columns

$wb = $excel.ActiveWorkbook

# ---- 1. Template metadata sheet: bump Version back down to 1.0.4 ----
$wsMeta = $wb.Worksheets.Item("isa_template")
$wsMeta.Range("B4").Value = "1.0.4"

# ---- 2. Annotation table sheet ----
$ws = $wb.Worksheets.Item("New Table")
$tbl = $ws.ListObjects.Item(1)

# Grow the table from 19 to 23 columns (4 new blank columns appended at the
# end); we then rewrite every header/data cell into its correct final
# position so the 4 new columns end up where they belong (two 3-column
# blocks inserted mid-table) instead of at the tail.
$tbl.Resize($ws.Range("A1:W2"))

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "Input [Sample Name]"
$ws.Range("B1").Value = "Component [next generation sequencing instrument model]"
$ws.Range("C1").Value = "Term Source REF (DPBO:0000040)"
$ws.Range("D1").Value = "Term Accession Number (DPBO:0000040)"
$ws.Range("E1").Value = "Characteristic [Sequencing Library Name]"
$ws.Range("F1").Value = "Term Source REF (NCIT:C182058)"
$ws.Range("G1").Value = "Term Accession Number (NCIT:C182058)"
$ws.Range("H1").Value = "Characteristic [Sequencing Library Source Indicator]"
$ws.Range("I1").Value = "Term Source REF (NCIT:C175895)"
$ws.Range("J1").Value = "Term Accession Number (NCIT:C175895)"
$ws.Range("K1").Value = "Parameter [Library selection]"
$ws.Range("L1").Value = "Term Source REF (GENEPIO:0001940)"
$ws.Range("M1").Value = "Term Accession Number (GENEPIO:0001940)"
$ws.Range("N1").Value = "Parameter [library strategy]"
$ws.Range("O1").Value = "Term Source REF (GENEPIO:0001973)"
$ws.Range("P1").Value = "Term Accession Number (GENEPIO:0001973)"
$ws.Range("Q1").Value = "Parameter [library layout]"
$ws.Range("R1").Value = "Term Source REF (DPBO:0000015)"
$ws.Range("S1").Value = "Term Accession Number (DPBO:0000015)"
$ws.Range("T1").Value = "Characteristic [MD5 Checksum]"
$ws.Range("U1").Value = "Term Source REF (NCIT:C171276)"
$ws.Range("V1").Value = "Term Accession Number (NCIT:C171276)"
$ws.Range("W1").Value = "Output [Data]"

# ---- Data row (row 2) ----
# A2 is already blank and keeps its position - leave untouched.
# B2:D2 are unaffected by the shift - leave untouched.
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = "Genomic DNA"
$ws.Range("I2").Value = "NCIT"
$ws.Range("J2").Value = "https://bioregistry.io/NCIT:C95940"
$ws.Range("K2").Value = "PCR method"
$ws.Range("L2").Value = "GENEPIO"
$ws.Range("M2").Value = "http://purl.obolibrary.org/obo/GENEPIO_0001955"
$ws.Range("N2").Value = "Whole Genome Sequencing"
$ws.Range("O2").Value = "NCIT"
$ws.Range("P2").Value = "https://bioregistry.io/NCIT:C101294"
$ws.Range("Q2").Value = "single-end"
$ws.Range("R2").Value = "DPBO"
$ws.Range("S2").Value = "http://purl.org/nfdi4plants/ontology/dpbo/DPBO_0000086"
# T2:W2 are brand new cells from the resize and are already blank.
